$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated nucleotide frequency values per position (publication re-run).
# Rows: 2=A, 3=C, 4=G, 5=T ; Columns B:X = positions 1..23

# Row 2
$ws.Range("B2").Value = 0.00681302043906132
$ws.Range("C2").Value = 0.00151400454201363
$ws.Range("D2").Value = 0.00151400454201363
$ws.Range("E2").Value = 0.000757002271006813
$ws.Range("F2").Value = 0.000757002271006813
$ws.Range("G2").Value = 0.00151400454201363
$ws.Range("H2").Value = 0.0416351249053747
$ws.Range("I2").Value = 0.0355791067373202
$ws.Range("J2").Value = 0.0105980317940954
$ws.Range("K2").Value = 0.0227100681302044
$ws.Range("L2").Value = 0.012112036336109
$ws.Range("M2").Value = 0.00454201362604088
$ws.Range("N2").Value = 0.00529901589704769
$ws.Range("O2").Value = 0.994700984102952
$ws.Range("P2").Value = 0.000757002271006813
$ws.Range("U2").Value = 0.00151400454201363
$ws.Range("V2").Value = 0.934140802422407
$ws.Range("W2").Value = 0.0151400454201363
$ws.Range("X2").Value = 0.0060560181680545

# Row 3
$ws.Range("B3").Value = 0.987887963663891
$ws.Range("C3").Value = 0.00757002271006813
$ws.Range("D3").Value = 0.0060560181680545
$ws.Range("E3").Value = 0.00681302043906132
$ws.Range("F3").Value = 0.989401968205905
$ws.Range("G3").Value = 0.98107494322483
$ws.Range("H3").Value = 0.0060560181680545
$ws.Range("I3").Value = 0.945495836487509
$ws.Range("J3").Value = 0.0280090840272521
$ws.Range("K3").Value = 0.000757002271006813
$ws.Range("M3").Value = 0.000757002271006813
$ws.Range("P3").Value = 0.998485995457986
$ws.Range("R3").Value = 0.0060560181680545
$ws.Range("S3").Value = 0.00151400454201363
$ws.Range("T3").Value = 0.999242997728993
$ws.Range("U3").Value = 0.993186979560939
$ws.Range("W3").Value = 0.00908402725208176

# Row 4
$ws.Range("B4").Value = 0.000757002271006813
$ws.Range("C4").Value = 0.00529901589704769
$ws.Range("E4").Value = 0.00227100681302044
$ws.Range("F4").Value = 0.000757002271006813
$ws.Range("G4").Value = 0.00529901589704769
$ws.Range("H4").Value = 0.94776684330053
$ws.Range("I4").Value = 0.00302800908402725
$ws.Range("J4").Value = 0.000757002271006813
$ws.Range("K4").Value = 0.9666919000757
$ws.Range("L4").Value = 0.987130961392884
$ws.Range("M4").Value = 0.993943981831946
$ws.Range("N4").Value = 0.993943981831946
$ws.Range("O4").Value = 0.000757002271006813
$ws.Range("P4").Value = 0.000757002271006813
$ws.Range("U4").Value = 0.000757002271006813
$ws.Range("V4").Value = 0.0643451930355791
$ws.Range("W4").Value = 0.9666919000757
$ws.Range("X4").Value = 0.98107494322483

# Row 5
$ws.Range("B5").Value = 0.00378501135503407
$ws.Range("C5").Value = 0.985616956850871
$ws.Range("D5").Value = 0.992429977289932
$ws.Range("E5").Value = 0.990158970476911
$ws.Range("F5").Value = 0.00908402725208176
$ws.Range("G5").Value = 0.012112036336109
$ws.Range("H5").Value = 0.00454201362604088
$ws.Range("I5").Value = 0.0158970476911431
$ws.Range("J5").Value = 0.960635881907646
$ws.Range("K5").Value = 0.00908402725208176
$ws.Range("N5").Value = 0.000757002271006813
$ws.Range("O5").Value = 0.00454201362604088
$ws.Range("R5").Value = 0.993943981831946
$ws.Range("S5").Value = 0.998485995457986
$ws.Range("T5").Value = 0.000757002271006813
$ws.Range("U5").Value = 0.00454201362604088
$ws.Range("V5").Value = 0.00151400454201363
$ws.Range("W5").Value = 0.00908402725208176
$ws.Range("X5").Value = 0.0128690386071158
